# Added two new Mac-Addresses
# Appends 10 new data rows (147-156) to the "master-reg_center_machine_devic"
# sheet, following the same column layout as the existing rows:
#   A=regcntr_id, B=machine_id, C=device_id, D=lang_code, E=is_active,
#   F=cr_by, G=cr_dtimes, H=eff_dtimes
# The new rows use a different cr_by value ("superadmin") than the prior
# rows ("superadmin()"), which introduces a new shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(10001, 10030, 3000166),
    @(10001, 10030, 3000167),
    @(10001, 10030, 3000168),
    @(10001, 10030, 3000169),
    @(10001, 10030, 3000170),
    @(10001, 10031, 3000171),
    @(10001, 10031, 3000172),
    @(10001, 10031, 3000173),
    @(10001, 10031, 3000174),
    @(10001, 10031, 3000175)
)

$startRow = 147
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $regcntr_id = $newRows[$i][0]
    $machine_id = $newRows[$i][1]
    $device_id  = $newRows[$i][2]

    $ws.Cells.Item($r, 1).Value = $regcntr_id
    $ws.Cells.Item($r, 2).Value = $machine_id
    $ws.Cells.Item($r, 3).Value = $device_id
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Match the author's final view state: scrolled down with A148 selected.
$ws.Range("A148").Select()
